# MCfCL-05-Lecture/Book1.xlsx edit
# Turns the single-sheet "absolute reference" demo workbook into a 4-sheet
# lecture workbook: relative reference / absolute reference / mixed
# reference / F4-key demo sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: rename existing sheet to "相対参照" (relative reference) and
# replace its absolute-reference formulas with relative ones.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "相対参照"

$ws1.Range("B2").Value = 1
$ws1.Range("B3").Value = 2
$ws1.Range("B4").Value = 3
$ws1.Range("B5").Value = 4
$ws1.Range("B6").Value = 5

$ws1.Range("D1").Formula = "=B1"
$ws1.Range("C2").Formula = "=A2"
$ws1.Range("D2").Formula = "=B2"
$ws1.Range("E2").Formula = "=C2"
$ws1.Range("D3").Formula = "=B3"

$excel.ActiveWindow.DisplayFormulas = $true
$excel.ActiveWindow.Zoom = 200
$ws1.Range("B15").Select()

# ---------------------------------------------------------------------
# Sheet 2: new sheet "絶対参照" (absolute reference) - the original demo
# content (same numbers, same $B$2 absolute formulas).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "絶対参照"

$ws2.Range("B2").Value = 1
$ws2.Range("B3").Value = 2
$ws2.Range("B4").Value = 3
$ws2.Range("B5").Value = 4
$ws2.Range("B6").Value = 5

$ws2.Range("D1").Formula = "=`$B`$2"
$ws2.Range("C2").Formula = "=`$B`$2"
$ws2.Range("D2").Formula = "=`$B`$2"
$ws2.Range("E2").Formula = "=`$B`$2"
$ws2.Range("D3").Formula = "=`$B`$2"

$excel.ActiveWindow.DisplayFormulas = $true
$excel.ActiveWindow.Zoom = 200
$ws2.Range("A16").Select()

# ---------------------------------------------------------------------
# Sheet 3: new sheet "混合参照" (mixed reference).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "混合参照"

$ws3.Range("A1").Value = "A1"
$ws3.Range("A2").Value = "A2"
$ws3.Range("A3").Value = "A3"
$ws3.Range("A4").Value = "A4"
$ws3.Range("A5").Value = "A5"
$ws3.Range("B1").Value = "B1"
$ws3.Range("B2").Value = "B2"
$ws3.Range("B3").Value = "B3"
$ws3.Range("B4").Value = "B4"
$ws3.Range("B5").Value = "B5"
$ws3.Range("C1").Value = "C1"
$ws3.Range("C2").Value = "C2"
$ws3.Range("C3").Value = "C3"
$ws3.Range("C4").Value = "C4"
$ws3.Range("C5").Value = "C5"

$ws3.Range("D1").Formula = "=A`$1"
$ws3.Range("E1:F1").Formula = "=B`$1"
$ws3.Range("D2:D5").Formula = "=A`$1"
$ws3.Range("E2:E5").Formula = "=B`$1"
$ws3.Range("F2:F5").Formula = "=C`$1"
$ws3.Range("A6").Formula = "=`$A1"
$ws3.Range("B6:C6").Formula = "=`$A1"
$ws3.Range("A7:C10").Formula = "=`$A2"

$excel.ActiveWindow.DisplayFormulas = $true
$excel.ActiveWindow.Zoom = 200
$ws3.Range("D9").Select()

# ---------------------------------------------------------------------
# Sheet 4: new sheet "F4-key" demo.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "F4-key"

$ws4.Range("B2").Value = "参照元"
$ws4.Range("D2").Formula = "=B`$2"

$excel.ActiveWindow.DisplayFormulas = $true
$excel.ActiveWindow.Zoom = 200
$ws4.Range("D2").Select()
$ws4.Activate()
